# Daily auto push update: insert two new "2026/01/06" (Tuesday) entries
# into the time-ranking table on Sheet1, pushing the existing
# 2026/12/29 .. 2027/01/05 block down by two rows (577:618 -> 579:620).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new blank rows right before the current row 577.
$ws.Rows.Item(577).Insert()
$ws.Rows.Item(577).Insert()

# Row 577: 2026/01/06, 火, 18, 201
# Force the date column to be stored as literal text (matching the rest
# of the sheet, which uses plain strings rather than date serials), then
# strip the temporary text number format so the cell keeps the sheet's
# default (unstyled) appearance.
$ws.Cells.Item(577, 1).NumberFormat = "@"
$ws.Cells.Item(577, 1).Value = "2026/01/06"
$ws.Cells.Item(577, 1).ClearFormats()
$ws.Cells.Item(577, 2).Value = "火"
$ws.Cells.Item(577, 3).Value = 18
$ws.Cells.Item(577, 4).Value = 201

# Row 578: 2026/01/06, 火, 19, 201
$ws.Cells.Item(578, 1).NumberFormat = "@"
$ws.Cells.Item(578, 1).Value = "2026/01/06"
$ws.Cells.Item(578, 1).ClearFormats()
$ws.Cells.Item(578, 2).Value = "火"
$ws.Cells.Item(578, 3).Value = 19
$ws.Cells.Item(578, 4).Value = 201
